$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (avoid Excel auto-numeric conversion)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.170.72"
$ws.Range("E2").Value = "  -3.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.672.74"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.88"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.01"
$ws.Range("E6").Value = "  -7.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.671.47"
$ws.Range("E7").Value = "  -3.52%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("E11").Value = "  -5.38%  "
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.53"
$ws.Range("E13").Value = "  -6.11%  "
$ws.Range("E14").Value = "  -6.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.286.08"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.674.83"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.249.17"
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("E19").Value = "  -6.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.87"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.81"
$ws.Range("E21").Value = "  -4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.05"
$ws.Range("E22").Value = "  -6.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.22"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  -7.46%  "
$ws.Range("E26").Value = "  -4.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.12"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -5.94%  "
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.58"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.810.81"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.612.94"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("E36").Value = "  -7.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.986"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.74"
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("E40").Value = "  -7.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.322"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "440.87"
$ws.Range("E42").Value = "  -8.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.54"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  -7.01%  "
$ws.Range("E45").Value = "  -8.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.30"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "141.75"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "39.63"
$ws.Range("E49").Value = "  -10.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.749.59"
$ws.Range("E50").Value = "  -6.82%  "
$ws.Range("E51").Value = "  -5.02%  "

Write-Output "Applied cryptos update"
